$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'91.240.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.66%  "

$ws.Range("D3").Value = "'3.177.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.85%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'216.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.66%  "

$ws.Range("D6").Value = "'628.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.74%  "

$ws.Range("D7").Value = "'1.16"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +32.86%  "

$ws.Range("E8").Value = "  +2.96%  "

$ws.Range("E9").Value = "  -0.06%  "

$ws.Range("D10").Value = "'3.174.79"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.86%  "

$ws.Range("D11").Value = "'0.763"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +14.86%  "

$ws.Range("E12").Value = "  +7.95%  "

$ws.Range("D13").Value = "'0.0000246"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.96%  "

$ws.Range("D14").Value = "'5.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.08%  "

$ws.Range("D15").Value = "'35.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +9.16%  "

$ws.Range("D16").Value = "'90.912.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.73%  "

$ws.Range("D17").Value = "'3.754.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.58%  "

$ws.Range("D18").Value = "'3.166.37"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.16%  "

$ws.Range("D19").Value = "'3.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +12.33%  "

$ws.Range("D20").Value = "'14.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +8.97%  "

$ws.Range("D21").Value = "'475.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +12.36%  "

$ws.Range("D22").Value = "'0.0000213"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.99%  "

$ws.Range("D23").Value = "'9.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +11.85%  "

$ws.Range("D24").Value = "'5.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.61%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'96.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +17.04%  "

$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D26").Value = "'5.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +10.69%  "

$ws.Range("D27").Value = "'12.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.14%  "

$ws.Range("D28").Value = "'3.335.93"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.47%  "

$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("D30").Value = "'9.41"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +13.34%  "

$ws.Range("E31").Value = "  +1.21%  "

$ws.Range("E32").Value = "  +0.11%  "

$ws.Range("D33").Value = "'27.65"
$ws.Range("D33").Style = "Normal"

$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").Value = "'526.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.41%  "

$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "'0.190"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +41.96%  "

$ws.Range("E36").Value = "  +7.99%  "

$ws.Range("D37").Value = "'3.66"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.27%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.145"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.59%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").Value = "'7.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.76%  "

$ws.Range("D40").Value = "'1.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.44%  "

$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").Value = "'0.0871"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +27.03%  "

$ws.Range("B42").Value = "WhiteBITCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D42").Value = "'22.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.10%  "

$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").Value = "'0.421"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +17.33%  "

$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.07%  "

$ws.Range("E45").Value = "  +9.30%  "

$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "'0.715"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +21.58%  "

$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "'152.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.92%  "

$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").Value = "'4.70"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +12.18%  "

$ws.Range("D50").Value = "'1.36"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +12.22%  "

$ws.Range("D51").Value = "'45.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.35%  "
